# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" row at the top of the "总计" (summary) sheet,
#    pushing the existing 2022-Q3 / 2022-Q2 / 2022-Q1 rows down.
# 2. Insert a brand-new "2022-Q4" worksheet (with the quarter's fund
#    holdings) right after "总计", before the existing "2022-Q3" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update the "总计" summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push existing data rows (2-4) down to (3-5) and open up a fresh row 2.
$summary.Rows.Item(2).Insert()

# The inserted row copied formatting from the row above (the header row);
# strip that so B2:D2 end up unstyled like their siblings in rows 3-5.
$summary.Range("B2:D2").ClearFormats()

# Give A2 the same style as the other index cells in column A (s="2").
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# New 2022-Q4 summary row.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.35

# Renumber the index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q4" worksheet before "2022-Q3"
# ---------------------------------------------------------------------
$existingQ3 = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($existingQ3)
$q4.Name = "2022-Q4"

# Re-use the header/index-column formatting from the (about to shift)
# "2022-Q3" sheet so the new sheet's styles line up with its siblings.
$q3 = $wb.Worksheets.Item(3)
$q3.Range("A1:H1").Copy()
$q4.Range("A1:H1").PasteSpecial(-4122)
$q3.Range("A2").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# These columns hold text that looks numeric (fund codes with leading
# zeros, fixed-decimal percentages) - force text format so Excel doesn't
# silently coerce them into numbers and drop the formatting.
$q4.Range("B2:B5").NumberFormat = "@"
$q4.Range("D2:D5").NumberFormat = "@"
$q4.Range("E2:E5").NumberFormat = "@"
$q4.Range("F2:F5").NumberFormat = "@"
$q4.Range("G2:G5").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "000739"
$q4.Range("C2").Value = "平安新鑫先锋混合A"
$q4.Range("D2").Value = "7.70"
$q4.Range("E2").Value = "86.57"
$q4.Range("F2").Value = "2.35"
$q4.Range("G2").Value = "0.1810"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "001515"
$q4.Range("C3").Value = "平安新鑫先锋混合C"
$q4.Range("D3").Value = "4.66"
$q4.Range("E3").Value = "86.57"
$q4.Range("F3").Value = "2.35"
$q4.Range("G3").Value = "0.1095"
$q4.Range("H3").Value = 9

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "011807"
$q4.Range("C4").Value = "平安研究精选混合A"
$q4.Range("D4").Value = "1.23"
$q4.Range("E4").Value = "86.30"
$q4.Range("F4").Value = "2.62"
$q4.Range("G4").Value = "0.0322"
$q4.Range("H4").Value = 8

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "011808"
$q4.Range("C5").Value = "平安研究精选混合C"
$q4.Range("D5").Value = "1.05"
$q4.Range("E5").Value = "86.30"
$q4.Range("F5").Value = "2.62"
$q4.Range("G5").Value = "0.0275"
$q4.Range("H5").Value = 8
